$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "videoRec" column (I) alongside the existing condition/recording
# columns (E:H) -- flags which trials had video recording enabled.
$ws.Range("I1").Value = "videoRec"

$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("I11").Value = 0

# Update the active selection to the new header cell; this also clears the
# stale scrolled-down view (topLeftCell) in favor of the default top-left
# view.
$ws.Range("I1").Select()
